$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 6
    3  = 6
    4  = 5
    5  = 7
    6  = 10
    7  = 5
    8  = 10
    9  = 4
    10 = 5
    11 = 6
    12 = 5
    13 = 4
    14 = 3
    15 = 7
    16 = 1
    17 = 4
    18 = 13
    19 = 7
    20 = 6
    21 = 9
    22 = 6
    23 = 5
    24 = 7
    25 = 5
    26 = 7
    27 = 9
    28 = 12
    29 = 8
    30 = 3
    31 = 11
    32 = 7
    33 = 5
    34 = 9
    35 = 5
    36 = 3
    37 = 5
    38 = 4
    39 = 2
    40 = 2
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
